$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "89.432.35"
$ws.Range("E2").Value = "  -1.45%  "
$ws.Range("D3").Value = "3.143.08"
$ws.Range("E3").Value = "  -1.11%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.52"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +1.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "637.08"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +3.58%  "
$ws.Range("E7").Value = "  +1.55%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.797"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  +16.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").Value = "3.137.32"
$ws.Range("E10").Value = "  -1.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.564"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -1.66%  "
$ws.Range("E12").Value = "  +1.46%  "
$ws.Range("E13").Value = "  -0.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.35"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +2.67%  "
$ws.Range("D15").Value = "89.196.37"
$ws.Range("E15").Value = "  -1.42%  "
$ws.Range("D16").Value = "3.708.14"
$ws.Range("E16").Value = "  -1.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "32.26"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -1.41%  "
$ws.Range("D18").Value = "3.138.15"
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.42"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +5.47%  "
$ws.Range("E20").Value = "  +20.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.24"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -1.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "423.97"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -2.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.47"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -0.63%  "
$ws.Range("E24").Value = "  -3.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.43"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +5.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "82.56"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +10.41%  "
$ws.Range("E27").Value = "  -2.55%  "
$ws.Range("D28").Value = "3.297.80"
$ws.Range("E28").Value = "  -2.50%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("E31").Value = "  -7.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.00"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -4.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.19"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -3.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "505.13"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.148"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +17.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.96"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +1.07%  "
$ws.Range("E37").Value = "  +3.16%  "
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.37"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +2.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.27"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -0.16%  "
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.367"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -2.03%  "
$ws.Range("E44").Value = "  -2.69%  "
$ws.Range("E45").Value = "  +8.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "146.20"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.70"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -2.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "163.87"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -5.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0656"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +11.73%  "
$ws.Range("E50").Value = "  +2.54%  "
$ws.Range("B51").Value = "ImmutableX"
$ws.Range("C51").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.19"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -3.50%  "
